$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.834.65"
$ws.Range("E2").Value = "  +0.37%  "
$ws.Range("D3").Value = "2.356.35"
$ws.Range("E3").Value = "  -0.45%  "
$ws.Range("E4").Value = "  +0.11%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "314.95"
$ws.Range("E5").Value = "  -4.64%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "108.82"
$ws.Range("E6").Value = "  +8.32%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.637"
$ws.Range("E7").Value = "  +0.13%  "
$ws.Range("E8").Value = "  -0.08%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.615"
$ws.Range("E9").Value = "  -0.34%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "40.62"
$ws.Range("E10").Value = "  +0.97%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0927"
$ws.Range("E11").Value = "  +0.63%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "8.55"
$ws.Range("E12").Value = "  +1.91%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.00"
$ws.Range("E13").Value = "  -1.42%  "
$ws.Range("E14").Value = "  +1.08%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "15.90"
$ws.Range("E15").Value = "  -2.11%  "
$ws.Range("D16").Value = "2.712.69"
$ws.Range("E16").Value = "  -0.42%  "
$ws.Range("D17").Value = "2.422.43"
$ws.Range("E17").Value = "  +2.20%  "
$ws.Range("D18").Value = "42.850.92"
$ws.Range("E18").Value = "  +0.40%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.62"
$ws.Range("E19").Value = "  -1.80%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.0000106"
$ws.Range("E20").Value = "  +0.03%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "76.59"
$ws.Range("E21").Value = "  +1.85%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "3.58"
$ws.Range("E22").Value = "  -3.88%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "270.03"
$ws.Range("E23").Value = "  -1.77%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.33"
$ws.Range("E24").Value = "  +0.40%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.51"
$ws.Range("E25").Value = "  -2.08%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "11.38"
$ws.Range("E27").Value = "  -0.81%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "23.41"
$ws.Range("E28").Value = "  -1.92%  "
$ws.Range("E29").Value = "  +2.15%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "36.62"
$ws.Range("E30").Value = "  +3.41%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "167.16"
$ws.Range("E31").Value = "  -3.74%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.0908"
$ws.Range("E32").Value = "  +0.70%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "6.12"
$ws.Range("E33").Value = "  +4.37%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.92"
$ws.Range("E34").Value = "  -6.21%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.121"
$ws.Range("E35").Value = "  +15.68%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.131"
$ws.Range("E36").Value = "  -0.40%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "4.69"
$ws.Range("E37").Value = "  +1.38%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0359"
$ws.Range("E38").Value = "  +0.07%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.81"
$ws.Range("E39").Value = "  -2.73%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.66"
$ws.Range("E40").Value = "  -8.27%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "105.51"
$ws.Range("E41").Value = "  +16.00%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.50"
$ws.Range("E42").Value = "  -0.50%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.237"
$ws.Range("E43").Value = "  +4.08%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "71.52"
$ws.Range("E45").Value = "  +0.12%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "12.52"
$ws.Range("E46").Value = "  +4.79%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "114.35"
$ws.Range("E47").Value = "  -1.22%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "79.81"
$ws.Range("E48").Value = "  +17.18%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "5.53"
$ws.Range("E49").Value = "  +1.11%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "9.08"
$ws.Range("E50").Value = "  +1.01%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.28"
$ws.Range("E51").Value = "  +1.67%  "
